$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the fill style (s="3") previously applied to D7:D12 - revert to default style
$ws.Range("D7:D12").Style = "Normal"

# Add new row 20 data - sc16 test case
$ws.Range("A20").Value = "sc16"
$ws.Range("B20").Value = 2
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = "A sublimit and a restriction on one of two policies"
$ws.Range("F20").Value = "complete"
$ws.Range("G20").Value = "yes"
$ws.Range("H20").Value = "done"

# Update selection/view to match the target state
$ws.Range("F20").Select()
